$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A, shifting column B (and its contents/formatting) into column A
$ws.Range("A1").EntireColumn.Delete()
